$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell values (rows 2-5). The order of entry matters for how
# Excel builds the shared-strings table (it appends new unique strings
# in first-seen order), so fill column C (DevAddr) first, then column B
# (Name), then column A (DevEUI), and finally columns E/F (the shared
# key value) - matching how the source workbook was authored.

$ws.Range("C2").Value = "BADDAD01"
$ws.Range("C3").Value = "BADDAD02"
$ws.Range("C4").Value = "BADDAD03"
$ws.Range("C5").Value = "BADDAD04"

$ws.Range("B2").Value = "Workshop-Device-01"
$ws.Range("B3").Value = "Workshop-Device-02"
$ws.Range("B4").Value = "Workshop-Device-03"
$ws.Range("B5").Value = "Workshop-Device-04"

$ws.Range("A2").Value = "104A15C001FFFF01"
$ws.Range("A3").Value = "104A15C001FFFF02"
$ws.Range("A4").Value = "104A15C001FFFF03"
$ws.Range("A5").Value = "104A15C001FFFF04"

$ws.Range("E2").Value = "104A15C001104A15C001104A15C001FF"
$ws.Range("F2").Value = "104A15C001104A15C001104A15C001FF"
$ws.Range("E3").Value = "104A15C001104A15C001104A15C001FF"
$ws.Range("F3").Value = "104A15C001104A15C001104A15C001FF"
$ws.Range("E4").Value = "104A15C001104A15C001104A15C001FF"
$ws.Range("F4").Value = "104A15C001104A15C001104A15C001FF"
$ws.Range("E5").Value = "104A15C001104A15C001104A15C001FF"
$ws.Range("F5").Value = "104A15C001104A15C001104A15C001FF"

# Column width adjustments (values chosen so the stored OOXML column
# width - after this runtime's internal rounding - lands on, or as close
# as possible to, the target widths from the authored workbook)
$ws.Range("B:B").ColumnWidth = 30.5
$ws.Range("C:C").ColumnWidth = 12
$ws.Range("D:D").ColumnWidth = 14.666666666666666
$ws.Range("E:F").ColumnWidth = 49.333333333333336
$ws.Range("G:G").ColumnWidth = 18.666666666666668
$ws.Range("H:H").ColumnWidth = 39.666666666666664
$ws.Range("I:I").ColumnWidth = 13.666666666666666

# Selection
$ws.Range("B20").Select()
